$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing existing rows 16-75 down to 17-76.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new market-record data. Constant columns
# (A,B,C,E,F,G,H,I,N,Q,R) repeat the same template values used by every
# other row in this sheet.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44701
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100112035
$ws.Cells.Item(16, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 30
$ws.Cells.Item(16, 11).Value = 28000
$ws.Cells.Item(16, 12).Value = 30000
$ws.Cells.Item(16, 13).Value = 29333
$ws.Cells.Item(16, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 2933
$ws.Cells.Item(16, 17).Value = 10
$ws.Cells.Item(16, 18).Value = "Hortaliza"
